$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(902, 1, 0, 0, 0, 0),
    @(201, 9, 30, 15, 45, 30),
    @(301, 6, 45, 30, 60, 45),
    @(501, 9, 52, 30, 75, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(101, 9, 30, 15, 60, 15),
    @(901, 16, 15, 45, 60, 60),
    @(701, 3, 90, 45, 97, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(1001, 18, 30, 75, 60, 72),
    @(401, 9, 48, 67, 75, 45),
    @(601, 9, 60, 67, 60, 42),
    @(801, 3, 67, 65, 52, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(1, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

